# Append six new daily COVID overview rows (2021-09-21 .. 2021-09-26)
# to the bottom of the "covid_totals" sheet, growing the used range
# from A1:H405 to A1:H411.
#
# Columns: A=date, B=areaType, C=areaCode, D=areaName,
#          E=cumCasesByPublishDate, F=newCasesByPublishDate,
#          G=newDeaths28DaysByPublishDate, H=cumDeaths28DaysByPublishDate

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "date" column has always held plain text (e.g. "2021-09-20"), not
# real Excel date serials. Pre-format the new cells as Text so the
# COM layer doesn't auto-coerce the "YYYY-MM-DD" strings into dates.
$ws.Range("A406:A411").NumberFormat = "@"

$newRows = @(
    @("2021-09-21", "overview", "K02000001", "United Kingdom", 7496543, 31564, 203, 135455),
    @("2021-09-22", "overview", "K02000001", "United Kingdom", 7530103, 34460, 166, 135621),
    @("2021-09-23", "overview", "K02000001", "United Kingdom", 7565867, 36710, 182, 135803),
    @("2021-09-24", "overview", "K02000001", "United Kingdom", 7601487, 35623, 180, 135983),
    @("2021-09-25", "overview", "K02000001", "United Kingdom", 7631233, 31348, $null, $null),
    @("2021-09-26", "overview", "K02000001", "United Kingdom", 7664230, 32417, 58, 136168)
)

$row = 406
foreach ($rec in $newRows) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]

    # Row 410 (2021-09-25) has no G/H values in the source data, so
    # leave those cells untouched rather than writing blanks.
    if ($rec[6] -ne $null) {
        $ws.Cells.Item($row, 7).Value = $rec[6]
    }
    if ($rec[7] -ne $null) {
        $ws.Cells.Item($row, 8).Value = $rec[7]
    }

    $row++
}
